$d = $word.ActiveDocument


# --- Hunk 1: prepend "Need to consider " + relocate _GoBack bookmark ---
$r1 = $d.Content
$r1.Find.Execute("Find an alert", $true, $false, $false, $false, $false, $true, 1, $false, "Need to consider Find an alert", 2)

$r1b = $d.Content
$r1b.Find.Execute("Find an alert", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1c = $r1b.Duplicate
$r1c.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r1c)


# --- Hunk 2: new bullet before "Uninstall existing service" ---
$r2 = $d.Content
$r2.Find.Execute("Uninstall existing service", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2p = $r2.Duplicate
$r2p.Collapse(1)
$r2p.InsertParagraphBefore()
$r2xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Stop the existing service: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PushNotificationsForSunriverWebApp</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2p.InsertXML($r2xml)


# --- Hunk 3: consolidate "Change the settings in the config file" ---
$r3 = $d.Content
$r3.Find.Execute("Change the settings in the", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3para = $r3.Paragraphs(1).Range
$r3start = $r3para.Start
$r3del = $d.Range($r3start, $r3para.End - 1)
$r3del.Delete()
$r3ins = $d.Range($r3start, $r3start)
$r3xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"><w:body><w:p w:rsidR="009F6EE0" w:rsidRDefault="009F6EE0" w:rsidP="00246ABE"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Change the settings in the config file</w:t></w:r><w:r><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r3ins.InsertXML($r3xml)


# --- Hunk 4: consolidate "Go to the config file and change..." ---
$r4 = $d.Content
$r4.Find.Execute("Go to the", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4para = $r4.Paragraphs(1).Range
$r4start = $r4para.Start
$r4del = $d.Range($r4start, $r4para.End - 1)
$r4del.Delete()
$r4ins = $d.Range($r4start, $r4start)
$r4xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"><w:body><w:p w:rsidR="008F5621" w:rsidRDefault="00B074B9" w:rsidP="009F6EE0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Go to the config file and c</w:t></w:r><w:r><w:t>hange the location of the web service, per the original instructions</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r4ins.InsertXML($r4xml)
